# Backup QR Scanner data - 2025-12-29T09:10:29.461Z - Cache Bust: 1766999429461
#
# 1. Rename the worksheet from "Session" to "Neurology"
# 2. Append the latest scanner log rows (76-78) to the log sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet/session tab to reflect the subject being logged ---
$ws.Name = "Neurology"

# --- Append newly scanned rows ---
# Column A holds numeric-looking Student IDs that must stay as TEXT (the
# original log keeps them as text so leading data / formatting is never
# lost). Using the classic apostrophe text-prefix reproduces exactly what
# a user (or the scanner app) does when forcing Excel to store a number
# as text, which is what triggers the "Number Stored as Text" indicator
# that the rest of the sheet already has ignored.
$newRows = @(
    @{ Row = 76; A = "201888"; B = "Neurology"; C = "29/12/2025"; D = "11:00:55"; E = "Scan"; F = "emp17.farah.a.youssef@gmail.com" },
    @{ Row = 77; A = "201243"; B = "Neurology"; C = "29/12/2025"; D = "11:06:16"; E = "Scan"; F = "emp17.farah.a.youssef@gmail.com" },
    @{ Row = 78; A = "201479"; B = "Neurology"; C = "29/12/2025"; D = "11:10:24"; E = "Scan"; F = "emp17.farah.a.youssef@gmail.com" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row
    $ws.Range("A$rowNum").Value = "'" + $r.A
    $ws.Range("B$rowNum").Value = $r.B
    $ws.Range("C$rowNum").Value = $r.C
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("E$rowNum").Value = $r.E
    $ws.Range("F$rowNum").Value = $r.F
}

$wb.Save()
